# Add the new "excel index" fields to the Booklet table (Tabelle5):
# AssessmentType, Description, Disclaimer, Duration, EscoOccupationId,
# EscoSkills, Publisher, Title. Only AssessmentType gets a default value
# of 0 for the existing data rows, matching the source system's export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Copy the formatting of the last existing header cell so new headers
# pick up the same style as the rest of the header row.
$lastHeader = $ws.Range("BB1")
$lastHeader.Copy()

$names = @("AssessmentType", "Description", "Disclaimer", "Duration", "EscoOccupationId", "EscoSkills", "Publisher", "Title")

foreach ($name in $names) {
    $col = $lo.ListColumns.Add()
    $hdr = $col.Range.Item(1)
    $hdr.Value = $name
    $hdr.PasteSpecial(-4122)  # xlPasteFormats
}

# Fill the new AssessmentType column's data rows with 0 (first new
# column only - matches the source data, which leaves the other new
# columns blank for now).
$lastRow = $lo.ListRows.Count + 1
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 55).Value = 0
}

# Put the selection where the author left it after adding the columns.
$ws.Range("BJ2").Select()
